# Update latest output (run 53)

$wb = $excel.ActiveWorkbook

# --- Schedule sheet: summary cost / unit cost update ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 772.4833740000001
$wsSchedule.Range("F2").Value = 12.77254255952381

# --- Detailed sheet: price / type updates ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Rows that flipped from forecast -> historical (with updated price)
$wsDetailed.Range("B13").Value = 73.20003
$wsDetailed.Range("C13").Value = "historical"

$wsDetailed.Range("C14").Value = "historical"

# Remaining price-only corrections
$wsDetailed.Range("B17").Value = 24.30882
$wsDetailed.Range("B18").Value = 25.52278
$wsDetailed.Range("B19").Value = 28.96238
$wsDetailed.Range("B23").Value = 35.88
$wsDetailed.Range("B24").Value = 35.88
$wsDetailed.Range("B25").Value = 36.06092
$wsDetailed.Range("B26").Value = 48.81498
$wsDetailed.Range("B29").Value = 36.06
$wsDetailed.Range("B30").Value = 36.06029
$wsDetailed.Range("B31").Value = 27.99749
$wsDetailed.Range("B32").Value = 30.21678
$wsDetailed.Range("B33").Value = 36.05863
$wsDetailed.Range("B34").Value = 33.26446
$wsDetailed.Range("B35").Value = 8.26928
$wsDetailed.Range("B36").Value = -3.07654
$wsDetailed.Range("B37").Value = -2.99582
$wsDetailed.Range("B38").Value = -2.91605
$wsDetailed.Range("B39").Value = -2.86176
$wsDetailed.Range("B40").Value = 3.42645
$wsDetailed.Range("B41").Value = 9.65137
$wsDetailed.Range("B42").Value = 9.71463
$wsDetailed.Range("B43").Value = 19.54531
$wsDetailed.Range("B44").Value = 8.333019999999999
$wsDetailed.Range("B45").Value = 6.5735
$wsDetailed.Range("B46").Value = 30.19777
$wsDetailed.Range("B49").Value = 47.60925
